$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 2: convert numeric-looking text cells to real numbers
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 8000
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 4

# Add row 3
$ws.Range("A3").Value = "Cool Kids Club"

# B3/C3 must stay text (not get auto-converted to numbers) but also must
# not pick up a lingering "Text" number-format style, so force text via
# NumberFormat, assign the value, then reset the style back to Normal.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "11111111"
$ws.Range("B3:C3").Style = "Normal"

$ws.Range("D3").Value = "Computer Science and Engineering"
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 1600
$ws.Range("G3").Value = 4.95
